$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 2, 3, 4 (data rows) had their values cyclically rotated:
#   new Row2 = old Row4
#   new Row3 = old Row2
#   new Row4 = old Row3
# Capture the original values first so the rotation is computed from a
# consistent snapshot (not partially-overwritten data).

$row2 = @{
    D = $ws.Range("D2").Value2
    I = $ws.Range("I2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    N = $ws.Range("N2").Value2
    P = $ws.Range("P2").Value2
    Q = $ws.Range("Q2").Value2
}

$row3 = @{
    D = $ws.Range("D3").Value2
    I = $ws.Range("I3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    P = $ws.Range("P3").Value2
    Q = $ws.Range("Q3").Value2
}

$row4 = @{
    D = $ws.Range("D4").Value2
    I = $ws.Range("I4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    P = $ws.Range("P4").Value2
    Q = $ws.Range("Q4").Value2
}

function Set-RowValues($targetRow, $values) {
    $ws.Range("D$targetRow").Value = $values.D
    $ws.Range("I$targetRow").Value = $values.I
    $ws.Range("J$targetRow").Value = $values.J
    $ws.Range("K$targetRow").Value = $values.K
    $ws.Range("L$targetRow").Value = $values.L
    $ws.Range("M$targetRow").Value = $values.M
    $ws.Range("N$targetRow").Value = $values.N
    $ws.Range("P$targetRow").Value = $values.P
    $ws.Range("Q$targetRow").Value = $values.Q
}

Set-RowValues 2 $row4
Set-RowValues 3 $row2
Set-RowValues 4 $row3
